$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A64").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("A65").Value = 65

$ws.Range("H58").Select()
